# Generate Report for Handoff
# Updates the "b.md" row to reflect that file has been handed off:
#  - Overview sheet: status columns + latest handoff date for b.md
#  - zh-cn / de-de sheets: status, new handoff file name (hyperlink text), new handoff datetime

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-25-20 00:25:50"

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-20 00:25:47"

# Hyperlinked cells can only be updated as a full set in this environment,
# so rebuild all of them, preserving the original targets/order and only
# changing the text/value for D3.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/a.md", "", "", ".md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb29afe583b30b9534ebc3e81d1eaefb6a5042c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2ab311bc8b8075614de409cdb4455b5f3b3e4297/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8d26c82ef51c2b3d6b106bad40a1c526e88e62d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/b.md", "", "", "b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/b.md", "", "", ".md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb29afe583b30b9534ebc3e81d1eaefb6a5042c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2ab311bc8b8075614de409cdb4455b5f3b3e4297/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8d26c82ef51c2b3d6b106bad40a1c526e88e62d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-20 00:25:50"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/a.md", "", "", ".md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56801cbe40e999aff98551f9fd2620125ad6824c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e7d35bb3553cefaadf5625974ed1bd016ff5e536/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ee10045b80065f74853c665811839ae0729168e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/b.md", "", "", "b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/b.md", "", "", ".md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56801cbe40e999aff98551f9fd2620125ad6824c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e7d35bb3553cefaadf5625974ed1bd016ff5e536/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ee10045b80065f74853c665811839ae0729168e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

Write-Host "Report regenerated for handoff."
